# "complete monthly and re-run daily"
#
# The underlying per-punch counts for ID sp99004 (陳衣玲, 進貨控場) were
# updated, and the downstream daily / TL / team-function aggregates were
# re-run to reflect the corrected counts.

$wb = $excel.ActiveWorkbook

$wsTeam = $wb.Worksheets.Item("team_df")
$wsDay  = $wb.Worksheets.Item("team_df_day")
$wsTl   = $wb.Worksheets.Item("productivity_tl")
$wsFunc = $wb.Worksheets.Item("productivity_team_function")

# --- team_df: per-shift rows for sp99004 (陳衣玲, 進貨控場) ---

# Row 6 (shift ending ~16:00)
$wsTeam.Range("S6").Value = 3
$wsTeam.Range("T6").Value = 32
$wsTeam.Range("U6").Value = 3 / 32

# Row 12 (shift ending ~18:58)
$wsTeam.Range("T12").Value = 25
$wsTeam.Range("U12").Value = 2 / 25

# --- team_df_day: daily aggregate for sp99004 / 進貨控場 ---
$wsDay.Range("F2").Value = 5
$wsDay.Range("G2").Value = 57
$wsDay.Range("H2").Value = 5 / 57

# --- productivity_tl: TL productivity score for sp99004 ---
$wsTl.Range("D2").Value = 5 / 57

# --- productivity_team_function: date_on_duty ratio for sp99004 ---
$wsFunc.Range("D2").Value = 5 / 57
